# Apply "Remove unused fields and add distribution_name to entitlements
# and deliveries tables" edit to the entitlements.xlsx workbook.

$wb = $excel.ActiveWorkbook

# --- "model" sheet: the entitlements field list -----------------------
$model = $wb.Worksheets.Item("model")

# Remove the unused "authorization_description" row entirely (old row 2).
$model.Rows.Item(2).Delete()

# Rename the old "authorization_name" field (now row 3) to
# "distribution_name".
$model.Range("B3").Value2 = "distribution_name"

# --- "settings" sheet: instance_name now references distribution_name -
$settings = $wb.Worksheets.Item("settings")
$settings.Range("B6").Value2 = "distribution_name"

# --- Selections / active sheet, matching the saved UI state -----------
[void]$model.Range("B3").Select()

$survey = $wb.Worksheets.Item("survey")
[void]$survey.Range("C13").Select()

[void]$settings.Activate()
[void]$settings.Range("C21").Select()
